$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format before assigning values, so that
# numeric-looking strings (e.g. "1.008") are NOT auto-converted to numbers.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "29.872.10"
$ws.Range("E2").Value = "  +1.62%  "
$ws.Range("D3").Value = "1.935.07"
$ws.Range("E3").Value = "  +1.30%  "
$ws.Range("D4").Value = "1.008"
$ws.Range("E4").Value = "  -0.40%  "
$ws.Range("D5").Value = "335.45"
$ws.Range("E5").Value = "  +3.06%  "
$ws.Range("D6").Value = "1.006"
$ws.Range("E6").Value = "  -0.36%  "
$ws.Range("D7").Value = "0.4829"
$ws.Range("E7").Value = "  +0.45%  "
$ws.Range("D8").Value = "0.4112"
$ws.Range("E8").Value = "  +1.14%  "
$ws.Range("D9").Value = "0.08157"
$ws.Range("E9").Value = "  -0.46%  "
$ws.Range("D10").Value = "1.014"
$ws.Range("D11").Value = "23.64"
$ws.Range("E11").Value = "  +1.15%  "
$ws.Range("D12").Value = "1.954.43"
$ws.Range("E12").Value = "  +1.87%  "
$ws.Range("D13").Value = "6.082"
$ws.Range("E13").Value = "  +0.98%  "
$ws.Range("D14").Value = "7.286"
$ws.Range("E14").Value = "  +1.20%  "
$ws.Range("D15").Value = "91.09"
$ws.Range("E15").Value = "  +0.25%  "
$ws.Range("D16").Value = "0.06856"
$ws.Range("E16").Value = "  +0.91%  "
$ws.Range("D17").Value = "1.008"
$ws.Range("E18").Value = "  -0.28%  "
$ws.Range("D19").Value = "17.79"
$ws.Range("E19").Value = "  +0.82%  "
$ws.Range("D20").Value = "1.006"
$ws.Range("E20").Value = "  -0.41%  "
$ws.Range("D21").Value = "29.882.00"
$ws.Range("E21").Value = "  +1.54%  "
$ws.Range("D22").Value = "5.628"
$ws.Range("E22").Value = "  +0.17%  "
$ws.Range("D23").Value = "11.86"
$ws.Range("D25").Value = "2.202.90"
$ws.Range("E25").Value = "  +2.46%  "
$ws.Range("D26").Value = "6.688"
$ws.Range("E26").Value = "  +1.31%  "
$ws.Range("D27").Value = "156.68"
$ws.Range("E27").Value = "  -0.16%  "
$ws.Range("D28").Value = "20.05"
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("D29").Value = "2.094"
$ws.Range("E29").Value = "  -0.28%  "
$ws.Range("D30").Value = "121.24"
$ws.Range("E30").Value = "  +1.10%  "
$ws.Range("E31").Value = "  -1.15%  "
$ws.Range("D32").Value = "0.09623"
$ws.Range("E32").Value = "  +0.76%  "
$ws.Range("D33").Value = "5.574"
$ws.Range("E33").Value = "  -0.49%  "
$ws.Range("D34").Value = "1.416"
$ws.Range("E34").Value = "  +3.93%  "
$ws.Range("D35").Value = "3.547"
$ws.Range("E35").Value = "  -0.31%  "
$ws.Range("D36").Value = "0.06566"
$ws.Range("E36").Value = "  +7.45%  "
$ws.Range("D37").Value = "0.02283"
$ws.Range("E37").Value = "  +0.21%  "
$ws.Range("D38").Value = "1.205"
$ws.Range("E38").Value = "  +2.29%  "
$ws.Range("D39").Value = "0.5963"
$ws.Range("E39").Value = "  +0.11%  "
$ws.Range("D40").Value = "7.964"
$ws.Range("E40").Value = "  -0.80%  "
$ws.Range("D41").Value = "10.72"
$ws.Range("E41").Value = "  -0.19%  "
$ws.Range("D42").Value = "0.1845"
$ws.Range("E42").Value = "  -0.19%  "
$ws.Range("D43").Value = "2.486"
$ws.Range("E43").Value = "  +3.36%  "
$ws.Range("E44").Value = "  +1.84%  "
$ws.Range("D45").Value = "12.34"
$ws.Range("E45").Value = "  -0.75%  "
$ws.Range("D46").Value = "0.07478"
$ws.Range("E46").Value = "  -1.48%  "
$ws.Range("D47").Value = "0.5563"
$ws.Range("E47").Value = "  +0.05%  "
$ws.Range("D48").Value = "1.991"
$ws.Range("E48").Value = "  +2.26%  "
$ws.Range("D49").Value = "117.14"
$ws.Range("E49").Value = "  +0.03%  "
$ws.Range("B50").Value = "MXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D50").Value = "2.422"
$ws.Range("E50").Value = "  +0.00%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "72.54"
$ws.Range("E51").Value = "  +0.39%  "

# Restore the default style on column D so no stray formatting remains
# (values are already committed as text; this only resets appearance).
$dRange.Style = "Normal"

